$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = "Andrea Riolfatti"
$ws.Range("B23").Value = "ELIA BATTISTI | U.S. Guarna"
$ws.Range("C23").Value = "Maickol Azocar | MAI UNA GIOIA"
$ws.Range("D23").Value = "Leonardo Viola | Shark Attack"
$ws.Range("E23").Value = "Stefano  Galvagni | Clitoriders"
$ws.Range("F23").Value = "Davide  Bazzano  | iMontagna"
